# "Generate Report for Handoff"
# b.md has now been handed off for localization (zh-cn / de-de); update the
# localization-status report to reflect the new handoff package.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1992f831aa4df72c37b812b49d09caf1be49fc59/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fce8d2533c4381d521616d21ffc54ce9da9fbbca/e2e/b.md."

# --- Overview sheet: b.md row (row 3) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-02 04:40:41"

# --- zh-cn sheet: b.md row (row 3) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-09-02 04:40:37"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.15

# --- de-de sheet: b.md row (row 3) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-09-02 04:40:41"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.15
